$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting for the new "2020" column (Q) from the existing "2019"
# column (P), cell by cell, then overwrite the values per the target data.

# Header year cell
$ws.Range("P3").Copy($ws.Range("Q3"))
$ws.Range("Q3").Value = 2020

# Data rows 4-12 all get the "-" placeholder (shared string), formatted like
# the rest of the data rows (style family used by P5:P12)
$ws.Range("P5").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = "-"

$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = "-"

$ws.Range("P6").Copy($ws.Range("Q6"))
$ws.Range("Q6").Value = "-"

$ws.Range("P7").Copy($ws.Range("Q7"))
$ws.Range("Q7").Value = "-"

$ws.Range("P8").Copy($ws.Range("Q8"))
$ws.Range("Q8").Value = "-"

$ws.Range("P9").Copy($ws.Range("Q9"))
$ws.Range("Q9").Value = "-"

$ws.Range("P10").Copy($ws.Range("Q10"))
$ws.Range("Q10").Value = "-"

$ws.Range("P11").Copy($ws.Range("Q11"))
$ws.Range("Q11").Value = "-"

$ws.Range("P12").Copy($ws.Range("Q12"))
$ws.Range("Q12").Value = "-"

# Row 13 uses the bottom-border style family (P13)
$ws.Range("P13").Copy($ws.Range("Q13"))
$ws.Range("Q13").Value = "-"

# Match the saved selection from the source workbook
$ws.Range("P17").Select()
